$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newTimestamp = "2025-03-02 01:27:28"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
